# poloniex history - "updated in 15 min"
# Finish off the XRP/USDT0000002 trade (row 41 buy + new row 42 sell) and
# open the next trade, XRP/USDT0000003 (new row 43 buy, still in progress).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1) Clone the formatting (number formats / wrap text / etc.) of row 41
#    into the two brand-new rows 42 and 43 before touching any values.
# ---------------------------------------------------------------------
$ws.Range("A41:L41").Copy()
$ws.Range("A42:L43").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 2) Finish off row 41: the XRP/USDT0000002 buy is now DONE, stamp the
#    finalized date/time and the fee that was charged.
# ---------------------------------------------------------------------
$ws.Cells.Item(41,8).Value = "DONE"
$ws.Cells.Item(41,9).Value = 42860.711898148147
$ws.Cells.Item(41,10).Value = "0.48000000 XRP (0.15%)"

# ---------------------------------------------------------------------
# 3) Row 42: the matching sell leg that closes out XRP/USDT0000002.
# ---------------------------------------------------------------------
$ws.Cells.Item(42,1).Value = 42860.441886574074

$ws.Cells.Item(42,2).Value = "            Sell"
$ws.Cells.Item(42,2).Characters(13,4).Font.Color = 255

$ws.Cells.Item(42,3).Value = "        XRP"

$ws.Cells.Item(42,4).NumberFormat = "@"
$ws.Cells.Item(42,4).Value = "              0.09865736`n"
$ws.Range("D41").Copy()
$ws.Cells.Item(42,4).PasteSpecial(-4122)

$ws.Cells.Item(42,5).Value = "          0.1USDT"
$ws.Cells.Item(42,6).Value = "         320 XRP"
$ws.Cells.Item(42,7).Value = " XRP/USDT0000002"
$ws.Cells.Item(42,8).Value = "DONE"
$ws.Cells.Item(42,9).Value = 42860.737060185187
$ws.Cells.Item(42,10).Value = "0.08000000 USDT (0.25%)"

$ws.Cells.Item(42,11).Value = "       ~17%"
$ws.Cells.Item(42,11).Characters(6,6).Font.Color = 5287936

$ws.Cells.Item(42,12).Value = "1 day"

# ---------------------------------------------------------------------
# 4) Row 43: the opening buy leg of the next trade, XRP/USDT0000003,
#    still in progress (no finalized date yet).
# ---------------------------------------------------------------------
$ws.Cells.Item(43,1).Value = 42860.737060185187

$ws.Cells.Item(43,2).Value = "            Buy"
$ws.Cells.Item(43,2).Characters(13,3).Font.Color = 5287936

$ws.Cells.Item(43,3).Value = "        XRP"

$ws.Cells.Item(43,4).NumberFormat = "@"
$ws.Cells.Item(43,4).Value = "              0.09750980`n"
$ws.Range("D41").Copy()
$ws.Cells.Item(43,4).PasteSpecial(-4122)

$ws.Cells.Item(43,5).Value = "          0.085USDT"
$ws.Cells.Item(43,6).Value = "         375 XRP"
$ws.Cells.Item(43,7).Value = " XRP/USDT0000003"
$ws.Cells.Item(43,8).Value = "IN PROGRESS"
$ws.Cells.Item(43,11).Value = "     "

# row 41's format-paste source spills a blank J/L placeholder onto row 43
# (row 41 itself has no L column, and J43 must stay untouched/empty) -
# drop those stray cells so row 43 only carries A-I and K, like the target.
$ws.Cells.Item(43,10).ClearContents()
$ws.Cells.Item(43,12).ClearContents()

# ---------------------------------------------------------------------
# 5) Row heights for the two new rows match the rest of the log, and the
#    view scrolls/selects where the author ended up after the edit.
# ---------------------------------------------------------------------
$ws.Rows.Item(42).RowHeight = 14.25
$ws.Rows.Item(43).RowHeight = 14.25

$ws.Application.Goto($ws.Range("A26"))
$ws.Range("B49").Select()
